# 5.17 Add History Image & Action
# Adds two new trailing columns ("History", "HistoryAction") to the header
# row of Sheet1, widens the new column R a bit, and leaves the selection
# where the author left it after the edit (R10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1 (these also register the two new shared strings
# "History" / "HistoryAction" and extend the sheet dimension/row spans
# automatically on save).
$ws.Range("Q1").Value = "History"
$ws.Range("R1").Value = "HistoryAction"

# Give the new column a bit more breathing room (matches the author's
# custom width for column R as closely as this engine's column-width
# rounding allows).
$ws.Columns.Item(18).ColumnWidth = 13.62

# Leave the selection on the newly added column, as in the authored edit.
$ws.Range("R10").Select()
